# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Corrects GP_RANK/other ranking columns and rewrites the BF "Date" column
# from the "4-10-2013-14" label format to ISO "2014-04-10".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = "AD2"; Value = 8 },
    @{ Cell = "BF2"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD3"; Value = 8 },
    @{ Cell = "BF3"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD4"; Value = 8 },
    @{ Cell = "AW4"; Value = 6 },
    @{ Cell = "BF4"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD5"; Value = 8 },
    @{ Cell = "AZ5"; Value = 2 },
    @{ Cell = "BF5"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD6"; Value = 8 },
    @{ Cell = "AF6"; Value = 11 },
    @{ Cell = "BF6"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD7"; Value = 1 },
    @{ Cell = "BA7"; Value = 25 },
    @{ Cell = "BF7"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "D8"; Value = 79 },
    @{ Cell = "F8"; Value = 31 },
    @{ Cell = "G8"; Value = 0.608 },
    @{ Cell = "K8"; Value = 0.474 },
    @{ Cell = "M8"; Value = 22.9 },
    @{ Cell = "W8"; Value = 8.6 },
    @{ Cell = "Y8"; Value = 3.6 },
    @{ Cell = "Z8"; Value = 20 },
    @{ Cell = "AC8"; Value = 2.6 },
    @{ Cell = "AF8"; Value = 9 },
    @{ Cell = "AG8"; Value = 9 },
    @{ Cell = "AW8"; Value = 5 },
    @{ Cell = "BF8"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "D9"; Value = 78 },
    @{ Cell = "E9"; Value = 34 },
    @{ Cell = "G9"; Value = 0.436 },
    @{ Cell = "J9"; Value = 85.7 },
    @{ Cell = "K9"; Value = 0.449 },
    @{ Cell = "M9"; Value = 23.6 },
    @{ Cell = "N9"; Value = 0.363 },
    @{ Cell = "O9"; Value = 18.9 },
    @{ Cell = "P9"; Value = 26.2 },
    @{ Cell = "Q9"; Value = 0.721 },
    @{ Cell = "R9"; Value = 12.2 },
    @{ Cell = "T9"; Value = 45.2 },
    @{ Cell = "X9"; Value = 5.7 },
    @{ Cell = "Y9"; Value = 5.7 },
    @{ Cell = "Z9"; Value = 23 },
    @{ Cell = "AB9"; Value = 104.4 },
    @{ Cell = "AC9"; Value = -2.2 },
    @{ Cell = "AD9"; Value = 8 },
    @{ Cell = "AE9"; Value = 19 },
    @{ Cell = "AH9"; Value = 26 },
    @{ Cell = "AK9"; Value = 16 },
    @{ Cell = "AT9"; Value = 5 },
    @{ Cell = "BF9"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD10"; Value = 1 },
    @{ Cell = "AT10"; Value = 3 },
    @{ Cell = "BF10"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "D11"; Value = 77 },
    @{ Cell = "F11"; Value = 29 },
    @{ Cell = "G11"; Value = 0.623 },
    @{ Cell = "I11"; Value = 39.1 },
    @{ Cell = "L11"; Value = 9.300000000000001 },
    @{ Cell = "N11"; Value = 0.378 },
    @{ Cell = "P11"; Value = 21.2 },
    @{ Cell = "Q11"; Value = 0.754 },
    @{ Cell = "T11"; Value = 45.3 },
    @{ Cell = "U11"; Value = 23.2 },
    @{ Cell = "AA11"; Value = 19.5 },
    @{ Cell = "AB11"; Value = 103.6 },
    @{ Cell = "AC11"; Value = 4.8 },
    @{ Cell = "AD11"; Value = 30 },
    @{ Cell = "AO11"; Value = 25 },
    @{ Cell = "AR11"; Value = 16 },
    @{ Cell = "AT11"; Value = 4 },
    @{ Cell = "AV11"; Value = 22 },
    @{ Cell = "AW11"; Value = 15 },
    @{ Cell = "BA11"; Value = 26 },
    @{ Cell = "BC11"; Value = 5 },
    @{ Cell = "BF11"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD12"; Value = 8 },
    @{ Cell = "AT12"; Value = 2 },
    @{ Cell = "BC12"; Value = 6 },
    @{ Cell = "BF12"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD13"; Value = 1 },
    @{ Cell = "BF13"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD14"; Value = 1 },
    @{ Cell = "AH14"; Value = 27 },
    @{ Cell = "AL14"; Value = 12 },
    @{ Cell = "BF14"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD15"; Value = 8 },
    @{ Cell = "BA15"; Value = 27 },
    @{ Cell = "BF15"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD16"; Value = 8 },
    @{ Cell = "AF16"; Value = 11 },
    @{ Cell = "BF16"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD17"; Value = 8 },
    @{ Cell = "BF17"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD18"; Value = 8 },
    @{ Cell = "BF18"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD19"; Value = 8 },
    @{ Cell = "AH19"; Value = 18 },
    @{ Cell = "AZ19"; Value = 3 },
    @{ Cell = "BF19"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD20"; Value = 8 },
    @{ Cell = "BF20"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD21"; Value = 8 },
    @{ Cell = "AK21"; Value = 17 },
    @{ Cell = "BF21"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD22"; Value = 8 },
    @{ Cell = "AH22"; Value = 18 },
    @{ Cell = "BF22"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD23"; Value = 8 },
    @{ Cell = "AO23"; Value = 26 },
    @{ Cell = "BF23"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD24"; Value = 8 },
    @{ Cell = "BF24"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD25"; Value = 8 },
    @{ Cell = "AG25"; Value = 10 },
    @{ Cell = "AT25"; Value = 12 },
    @{ Cell = "BF25"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD26"; Value = 1 },
    @{ Cell = "BF26"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD27"; Value = 1 },
    @{ Cell = "BF27"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "D28"; Value = 78 },
    @{ Cell = "E28"; Value = 60 },
    @{ Cell = "G28"; Value = 0.769 },
    @{ Cell = "I28"; Value = 40.6 },
    @{ Cell = "K28"; Value = 0.489 },
    @{ Cell = "L28"; Value = 8.4 },
    @{ Cell = "M28"; Value = 21.1 },
    @{ Cell = "N28"; Value = 0.398 },
    @{ Cell = "Q28"; Value = 0.784 },
    @{ Cell = "R28"; Value = 9.199999999999999 },
    @{ Cell = "V28"; Value = 14.7 },
    @{ Cell = "AA28"; Value = 19.5 },
    @{ Cell = "AC28"; Value = 8.1 },
    @{ Cell = "AD28"; Value = 8 },
    @{ Cell = "AH28"; Value = 24 },
    @{ Cell = "AL28"; Value = 13 },
    @{ Cell = "AT28"; Value = 13 },
    @{ Cell = "AZ28"; Value = 1 },
    @{ Cell = "BA28"; Value = 24 },
    @{ Cell = "BF28"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD29"; Value = 8 },
    @{ Cell = "AF29"; Value = 11 },
    @{ Cell = "BF29"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD30"; Value = 8 },
    @{ Cell = "AR30"; Value = 17 },
    @{ Cell = "BF30"; Value = "2014-04-10"; IsText = $true },
    @{ Cell = "AD31"; Value = 8 },
    @{ Cell = "AR31"; Value = 18 },
    @{ Cell = "BF31"; Value = "2014-04-10"; IsText = $true }
)

foreach ($edit in $edits) {
    $range = $ws.Range($edit.Cell)
    if ($edit.ContainsKey("IsText") -and $edit.IsText) {
        # Force text type so Excel does not auto-convert date-looking
        # strings (e.g. "2014-04-10") into a date serial number.
        $range.NumberFormat = "@"
        $range.Value = $edit.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $edit.Value
    }
}
